$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Appreciations")

# Insert a new row at position 3 (shifts existing rows 3-7 down to 4-8)
$ws.Rows(3).Insert()

# Populate the new row 3 with a new appreciation entry.
# Core value / Core value description mirror the existing entries,
# the appreciation description is the new long text, and the
# sender/receiver/designation columns mirror row 2 (Samnit Patil, Software Engineer,
# Sharyu Marwadi, Trainee).
$ws.Cells.Item(3, 1).Value = "corevalue 1 updated"
$ws.Cells.Item(3, 2).Value = "new desc 1"
$ws.Cells.Item(3, 3).Value = "Great teamwork dskafkdfgkjahda a dhsfkjadhkfha hdkjfhakdhfkg sadhfkahdsjkfhakj aksjdhfkjahdkjfhakjsd akdshfkjahdfk"
$ws.Cells.Item(3, 4).Value = "Samnit"
$ws.Cells.Item(3, 5).Value = "Patil"
$ws.Cells.Item(3, 6).Value = "Software Engineer"
$ws.Cells.Item(3, 7).Value = "Sharyu"
$ws.Cells.Item(3, 8).Value = "Marwadi"
$ws.Cells.Item(3, 9).Value = "Trainee"
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
